# Auto-update draw results: append the 2025-11-20 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row goes right after the current last data row.
$newRow = $ws.UsedRange.Rows.Count + 1

# The existing "Date"/"Phase" columns hold numeric-looking text (e.g. "2025-11-20",
# "251120") that must stay as literal text, not be auto-converted to a date serial
# or a number. Force Text format on the row before writing the values so every
# cell (like the rest of the sheet) is stored as text.
$rng = $ws.Range("A" + $newRow + ":E" + $newRow)
$rng.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-11-20"
$ws.Range("B" + $newRow).Value = "Pick 4"
$ws.Range("C" + $newRow).Value = "251120"
$ws.Range("D" + $newRow).Value = "1-0-4-0"
$ws.Range("E" + $newRow).Value = "2025-11-20T21:38:41.562+04:00"
